# Append 1 row (row 12) at the bottom of the data on Sheet1, matching the
# "Append 1 rows" commit: A=<blank>, B=محمود, C=23 (kept as text),
# D=النصر, E=الرحلة 2, F=C3, G=IDRF, H=<arabic timestamp>.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = ""
$ws.Range("B12").Value = "محمود"
# Leading apostrophe forces text storage (matches the existing sheet's
# numberStoredAsText usage for C2:C11) without altering the cell's
# number format.
$ws.Range("C12").Value = "'23"
$ws.Range("D12").Value = "النصر"
$ws.Range("E12").Value = "الرحلة 2"
$ws.Range("F12").Value = "C3"
$ws.Range("G12").Value = "IDRF"
$ws.Range("H12").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:٥٦:٥٦ م"
